$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 340, shifting existing rows 340-351 down to 341-352
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new weekly price record
$ws.Cells.Item(340, 1).Value = 7
$ws.Cells.Item(340, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(340, 3).Value = "Ñuble"
$ws.Cells.Item(340, 4).Value = 45075
$ws.Cells.Item(340, 5).Value = 16
$ws.Cells.Item(340, 6).Value = 100112017
$ws.Cells.Item(340, 7).Value = "Apio"
$ws.Cells.Item(340, 8).Value = "Americana (o)"
$ws.Cells.Item(340, 9).Value = "Primera"
$ws.Cells.Item(340, 10).Value = 80
$ws.Cells.Item(340, 11).Value = 6000
$ws.Cells.Item(340, 12).Value = 6000
$ws.Cells.Item(340, 13).Value = 6000
$ws.Cells.Item(340, 14).Value = "$/docena de matas"
$ws.Cells.Item(340, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(340, 16).Value = 1000
$ws.Cells.Item(340, 17).Value = 6
$ws.Cells.Item(340, 18).Value = "Hortaliza"
